$d = $word.ActiveDocument

# 1. Merge "Standard " + "Nondisclosure Agreement" bold runs into a single run.
$d.Content.Find.Execute("Standard Nondisclosure Agreement", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Standard Nondisclosure Agreement", 2) | Out-Null

# 2. Remove the "https://" hyperlink before {{{domain}}}, unwrapping the
#    hyperlink field but keeping the bookmark in place.
$h = $d.Hyperlinks.Item(1)
$h.Delete()

# Now that the hyperlink field is gone, "https://" is a plain run of text
# that can be found & removed with a normal Find/Replace.
$d.Content.Find.Execute("https://", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# Move the space that used to sit after "published at " so that it now sits
# before "{{{domain}}}." instead (since "https://" no longer separates them).
$d.Content.Find.Execute("published at ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "published at", 2) | Out-Null
$d.Content.Find.Execute("{{{domain}}}.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " {{{domain}}}.", 2) | Out-Null
